# apps/comparator/res/concept-mappings.xlsx
# "added OR capabilities to comparator"
#
# - New mapping rows 27/28 gain a variable_name of "respiratory_rate" (col D).
# - Fix a typo: "oxygenation_index_cacl" -> "oxygenation_index_calc" (D29:D30).
# - Rename the "catecholamine-therapy" concept to "vasopressor-therapy"
#   (code in col B, display name in col C) across rows 33:37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Respiratory Rate concept rows now carry the variable_name they map to.
$ws.Range("D27").Value = "respiratory_rate"
$ws.Range("D28").Value = "respiratory_rate"

# Typo fix in the variable_name for the ARDS oxygenation index rows.
$ws.Range("D29").Value = "oxygenation_index_calc"
$ws.Range("D30").Value = "oxygenation_index_calc"

# Catecholamine Therapy -> Vasopressor Therapy rename (code + display name).
$ws.Range("B33").Value = "vasopressor-therapy"
$ws.Range("C33").Value = "Vasopressor Therapy"

$ws.Range("B34").Value = "vasopressor-therapy"
$ws.Range("C34").Value = "Vasopressor Therapy"

$ws.Range("B35").Value = "vasopressor-therapy"
$ws.Range("C35").Value = "Vasopressor Therapy"

$ws.Range("B36").Value = "vasopressor-therapy"
$ws.Range("C36").Value = "Vasopressor Therapy"

$ws.Range("B37").Value = "vasopressor-therapy"
$ws.Range("C37").Value = "Vasopressor Therapy"

# Match the author's final cursor/scroll position from the edit session.
[void]$ws.Range("D30").Select()
